$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column widths (stored XML width = ColumnWidth + 0.8333333333333333) ---
$ws.Range("B1:C1").ColumnWidth = 7.166666666666667
$ws.Range("F1:G1").ColumnWidth = 7.166666666666667
$ws.Range("I1:M1").ColumnWidth = 7.166666666666667
$ws.Range("O1:Q1").ColumnWidth = 7.166666666666667
$ws.Range("T1").ColumnWidth = 8.166666666666666
$ws.Range("U1:W1").ColumnWidth = 7.166666666666667
$ws.Range("Z1:AD1").ColumnWidth = 7.166666666666667
$ws.Range("AH1").ColumnWidth = 7.166666666666667

# --- Update cell values for rows 2-5 ---
# Row 2
$ws.Range("A2").Value = 45049.50694444445
$ws.Range("B2").Value = 10.726
$ws.Range("C2").Value = 7.333
$ws.Range("D2").Value = 3.404
$ws.Range("E2").Value = 23.56
$ws.Range("F2").Value = 17.15
$ws.Range("G2").Value = 8.176
$ws.Range("H2").Value = 24.228
$ws.Range("I2").Value = 13.347
$ws.Range("J2").Value = 5.245
$ws.Range("K2").Value = 7.323
$ws.Range("L2").Value = 9.308999999999999
$ws.Range("M2").Value = 10.191
$ws.Range("N2").Value = 2.44
$ws.Range("O2").Value = 8.647
$ws.Range("P2").Value = 11.655
$ws.Range("Q2").Value = 7.955
$ws.Range("R2").Value = 2.648
$ws.Range("S2").Value = 1.093
$ws.Range("T2").Value = 124.223
$ws.Range("U2").Value = 23.834
$ws.Range("V2").Value = 7.982
$ws.Range("W2").Value = 14.964
$ws.Range("X2").Value = 8.048999999999999
$ws.Range("Y2").Value = 2.19
$ws.Range("Z2").Value = 13.597
$ws.Range("AA2").Value = 7.05
$ws.Range("AB2").Value = 6.629
$ws.Range("AC2").Value = 7.562
$ws.Range("AD2").Value = 9.942
$ws.Range("AE2").Value = 2.682
$ws.Range("AF2").Value = 21.628
$ws.Range("AG2").Value = 4.075
$ws.Range("AH2").Value = 9.978

# Row 3
$ws.Range("A3").Value = 45049.51388888889
$ws.Range("B3").Value = 20.441
$ws.Range("C3").Value = 15.086
$ws.Range("D3").Value = 1.894
$ws.Range("E3").Value = 44.851
$ws.Range("F3").Value = 35.853
$ws.Range("G3").Value = 15.964
$ws.Range("H3").Value = 59.534
$ws.Range("I3").Value = 24.983
$ws.Range("J3").Value = 10.967
$ws.Range("K3").Value = 15.859
$ws.Range("L3").Value = 17.959
$ws.Range("M3").Value = 19.221
$ws.Range("N3").Value = 4.972
$ws.Range("O3").Value = 16.166
$ws.Range("P3").Value = 22.781
$ws.Range("Q3").Value = 13.928
$ws.Range("R3").Value = 1.372
$ws.Range("S3").Value = 0.894
$ws.Range("T3").Value = 238.748
$ws.Range("U3").Value = 45.164
$ws.Range("V3").Value = 14.922
$ws.Range("W3").Value = 29.986
$ws.Range("X3").Value = 15.78
$ws.Range("Y3").Value = 2.529
$ws.Range("Z3").Value = 30.159
$ws.Range("AA3").Value = 13.181
$ws.Range("AB3").Value = 11.837
$ws.Range("AC3").Value = 13.862
$ws.Range("AD3").Value = 18.902
$ws.Range("AE3").Value = 1.136
$ws.Range("AF3").Value = 54.309
$ws.Range("AG3").Value = 8.233000000000001
$ws.Range("AH3").Value = 18.656

# Row 4
$ws.Range("A4").Value = 45049.52083333334
$ws.Range("B4").Value = 8.003
$ws.Range("C4").Value = 5.863
$ws.Range("D4").Value = 1.023
$ws.Range("E4").Value = 17.731
$ws.Range("F4").Value = 13.823
$ws.Range("G4").Value = 6.223
$ws.Range("H4").Value = 28.958
$ws.Range("I4").Value = 9.859999999999999
$ws.Range("J4").Value = 4.326
$ws.Range("K4").Value = 5.989
$ws.Range("L4").Value = 7.098
$ws.Range("M4").Value = 7.675
$ws.Range("N4").Value = 1.892
$ws.Range("O4").Value = 6.391
$ws.Range("P4").Value = 8.949999999999999
$ws.Range("Q4").Value = 5.681
$ws.Range("R4").Value = 0.866
$ws.Range("S4").Value = 0.402
$ws.Range("T4").Value = 89.983
$ws.Range("U4").Value = 18.051
$ws.Range("V4").Value = 5.899
$ws.Range("W4").Value = 11.825
$ws.Range("X4").Value = 6.173
$ws.Range("Y4").Value = 1.091
$ws.Range("Z4").Value = 13.935
$ws.Range("AA4").Value = 5.211
$ws.Range("AB4").Value = 4.768
$ws.Range("AC4").Value = 5.568
$ws.Range("AD4").Value = 7.469
$ws.Range("AE4").Value = 0.722
$ws.Range("AF4").Value = 26.675
$ws.Range("AG4").Value = 3.169
$ws.Range("AH4").Value = 7.377

# Row 5
$ws.Range("A5").Value = 45049.52777777778
$ws.Range("B5").Value = 23.41
$ws.Range("C5").Value = 17.49
$ws.Range("D5").Value = 1.33
$ws.Range("E5").Value = 51.19
$ws.Range("F5").Value = 41.69
$ws.Range("G5").Value = 18.37
$ws.Range("H5").Value = 69.01000000000001
$ws.Range("I5").Value = 28.48
$ws.Range("J5").Value = 12.72
$ws.Range("K5").Value = 18.63
$ws.Range("L5").Value = 20.53
$ws.Range("M5").Value = 21.86
$ws.Range("N5").Value = 5.78
$ws.Range("O5").Value = 18.42
$ws.Range("P5").Value = 26.16
$ws.Range("Q5").Value = 15.56
$ws.Range("R5").Value = 0.8
$ws.Range("S5").Value = 0.8
$ws.Range("T5").Value = 273.08
$ws.Range("U5").Value = 51.45
$ws.Range("V5").Value = 17
$ws.Range("W5").Value = 34.55
$ws.Range("X5").Value = 18.16
$ws.Range("Y5").Value = 2.6
$ws.Range("Z5").Value = 34.25
$ws.Range("AA5").Value = 15.02
$ws.Range("AB5").Value = 13.33
$ws.Range("AC5").Value = 15.65
$ws.Range("AD5").Value = 21.59
$ws.Range("AE5").Value = 0.53
$ws.Range("AF5").Value = 62.58
$ws.Range("AG5").Value = 9.550000000000001
$ws.Range("AH5").Value = 21.26

# --- Delete row 6 (dataset trimmed from 5 to 4 data rows) ---
$ws.Rows.Item(6).Delete()
